$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades durante la pasantia")

# Update existing dates for rows 93-100 (shift ~1 month later)
$ws.Cells.Item(93, 1).Value = 44805
$ws.Cells.Item(94, 1).Value = 44807
$ws.Cells.Item(95, 1).Value = 44810
$ws.Cells.Item(96, 1).Value = 44812
$ws.Cells.Item(97, 1).Value = 44817
$ws.Cells.Item(98, 1).Value = 44819
$ws.Cells.Item(99, 1).Value = 44824
$ws.Cells.Item(100, 1).Value = 44826

# Add two new activity rows (102, 103), copying formatting from row 101
[void]$ws.Range("A101:C101").Copy()
[void]$ws.Range("A102:C102").PasteSpecial(-4122)
[void]$ws.Range("A103:C103").PasteSpecial(-4122)

$ws.Cells.Item(102, 1).Value = 44833
$ws.Cells.Item(102, 2).Value = "Reunion con vero y profe"
$ws.Cells.Item(102, 3).Value = 2

$ws.Cells.Item(103, 1).Value = 44839
$ws.Cells.Item(103, 2).Value = "Revision de datos atipicos componentes"
$ws.Cells.Item(103, 3).Value = 2

[void]$ws.Cells.Item(103, 2).Select()
